# Update dayrealmob correct cell in time
# For each data row (2-27), update columns B..L (Unnamed: 0, Attesa, Risposte Eff.,
# Offerte, Abb sup. 14, Abb inf. 14, T.A., Livello di Servizio %, Short Call min 10,
# Cleared, % Cleared) with the corrected figures. Column A (fcst_name) is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: row number, then values for columns B,C,D,E,F,G,H,I,J,K,L
$data = @(
  @(2, 3903, 3, 75, 180, 4, 0, 1096, 40465, 1, 95, 14656),
  @(3, 10, 1, 352, 353, 0, 0, 0, 8617, 6, 0, 0),
  @(4, 271, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(5, 804, 0, 173, 173, 0, 0, 0, 8830, 7, 0, 0),
  @(6, 272, 0, 25, 25, 0, 0, 0, 8000, 1, 0, 0),
  @(7, 269, 0, 7, 7, 0, 0, 0, 8571, 0, 0, 0),
  @(8, 1, 10, 131, 154, 6, 7, 408, 4733, 0, 0, 0),
  @(9, 776, 0, 28, 28, 0, 0, 0, 9643, 0, 0, 0),
  @(10, 785, 0, 31, 34, 3, 0, 1364, 18882, 2, 0, 0),
  @(11, 1046, 0, 143, 154, 9, 2, 2362, 40293, 2, 0, 0),
  @(12, 866, 0, 18, 18, 0, 0, 0, 15256, 1, 0, 0),
  @(13, 264, 4, 146, 179, 28, 1, 4739, 5620, 6, 0, 0),
  @(14, 805, 0, 144, 145, 0, 1, 0, 9225, 5, 0, 0),
  @(15, 270, 0, 10, 10, 0, 0, 0, 8000, 0, 0, 0),
  @(16, 236, 2, 231, 257, 7, 1, 1902, 15036, 3, 17, 909),
  @(17, 72, 0, 12, 12, 0, 0, 0, 7500, 0, 0, 0),
  @(18, 7, 1, 246, 259, 1, 1, 39, 8082, 8, 10, 386),
  @(19, 788, 0, 93, 94, 1, 0, 127, 9307, 0, 0, 0),
  @(20, 1407, 0, 89, 113, 21, 3, 2059, 4336, 0, 0, 0),
  @(21, 507, 0, 38, 38, 0, 0, 0, 20238, 2, 0, 0),
  @(22, 196, 0, 1, 1, 0, 0, 0, 10000, 0, 0, 0),
  @(23, 814, 0, 8, 8, 0, 0, 0, 8750, 1, 0, 0),
  @(24, 1013, 0, 17, 17, 0, 0, 0, 9286, 0, 0, 0),
  @(25, 16, 1, 23, 24, 0, 0, 0, 7391, 0, 0, 0),
  @(26, 142, 0, 17, 19, 2, 0, 1053, 7059, 0, 0, 0),
  @(27, 193, 0, 41, 41, 0, 1, 0, 5366, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
